$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 32.142857
$ws.Range("I6").Value = 3.8333333
$ws.Range("K6").Value = 11.4999999
$ws.Range("M6").Value = 100.5000001
$ws.Range("H20").Value = 1466
$ws.Range("I20").Value = 1466
$ws.Range("K20").Value = 1466
$ws.Range("M20").Value = -1236
$ws.Range("H33").Value = 186.625
$ws.Range("I33").Value = 113.28571
$ws.Range("J33").Value = 700
$ws.Range("K33").Value = 113.28571
$ws.Range("L33").Value = 700
$ws.Range("M33").Value = 115.71429
$ws.Range("N33").Value = -1158
$ws.Range("H35").Value = 1466
$ws.Range("I35").Value = 1466
$ws.Range("K35").Value = 1466
$ws.Range("M35").Value = -1087
$ws.Range("H51").Value = 6374.6665
$ws.Range("I51").Value = 4600
$ws.Range("K51").Value = 4600
$ws.Range("M51").Value = -4116
$ws.Range("H69").Value = 35875.438
$ws.Range("J69").Value = 24071.215
$ws.Range("L69").Value = 72213.645
$ws.Range("N69").Value = -73961.645
$ws.Range("H70").Value = 1443.3334
$ws.Range("I70").Value = 990
$ws.Range("J70").Value = 1500
$ws.Range("K70").Value = 2970
$ws.Range("L70").Value = 4500
$ws.Range("M70").Value = -2700
$ws.Range("N70").Value = -5040
$ws.Range("H72").Value = 35875.438
$ws.Range("J72").Value = 24071.215
$ws.Range("L72").Value = 216640.935
$ws.Range("N72").Value = -225376.935
$ws.Range("H73").Value = 1443.3334
$ws.Range("I73").Value = 990
$ws.Range("J73").Value = 1500
$ws.Range("K73").Value = 2970
$ws.Range("L73").Value = 4500
$ws.Range("M73").Value = -2034
$ws.Range("N73").Value = -6372
$ws.Range("H86").Value = 8000
$ws.Range("I86").Value = 8000
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 8000
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = $null
$ws.Range("N86").Value = -6877
$ws.Range("H89").Value = 8000
$ws.Range("I89").Value = 8000
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 40000
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = $null
$ws.Range("N89").Value = -34384
$ws.Range("H137").Value = 2460
$ws.Range("I137").Value = 2433.3333
$ws.Range("K137").Value = 7299.999899999999
$ws.Range("M137").Value = -4749.999899999999
$ws.Range("H138").Value = 3163.8333
$ws.Range("I138").Value = 3474
$ws.Range("J138").Value = 2543.5
$ws.Range("K138").Value = 10422
$ws.Range("L138").Value = 7630.5
$ws.Range("M138").Value = -5282
$ws.Range("N138").Value = -17910.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3076.4814
$ws.Range("I32").Value = 3164.577
$ws.Range("K32").Value = 3164.577
$ws.Range("M32").Value = -2877.577
$ws.Range("H56").Value = 39900
$ws.Range("J56").Value = 39900
$ws.Range("L56").Value = 39900
$ws.Range("N56").Value = -41384
$ws.Range("H61").Value = 1350
$ws.Range("I61").Value = 1520
$ws.Range("K61").Value = 1520
$ws.Range("M61").Value = -1308
$ws.Range("H64").Value = 27000
$ws.Range("J64").Value = 27000
$ws.Range("L64").Value = 27000
$ws.Range("N64").Value = -27496
$ws.Range("H67").Value = 27000
$ws.Range("J67").Value = 27000
$ws.Range("L67").Value = 27000
$ws.Range("N67").Value = -28716
$ws.Range("H68").Value = 52000
$ws.Range("J68").Value = 52000
$ws.Range("L68").Value = 52000
$ws.Range("N68").Value = -53622
$ws.Range("H71").Value = 52000
$ws.Range("J71").Value = 52000
$ws.Range("L71").Value = 156000
$ws.Range("N71").Value = -164112
$ws.Range("H74").Value = 1103.5555
$ws.Range("I74").Value = 989.75
$ws.Range("K74").Value = 989.75
$ws.Range("M74").Value = -115.75
$ws.Range("H77").Value = 1103.5555
$ws.Range("I77").Value = 989.75
$ws.Range("K77").Value = 4948.75
$ws.Range("M77").Value = -580.75
$ws.Range("H110").Value = 6356.4
$ws.Range("I110").Value = 6356.4
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 6356.4
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = $null
$ws.Range("N110").Value = -4311.4
$ws.Range("H132").Value = 972.0769
$ws.Range("I132").Value = 968.8182
$ws.Range("K132").Value = 2906.4546
$ws.Range("M132").Value = -376.4546
$ws.Range("H136").Value = 1350
$ws.Range("I136").Value = 1520
$ws.Range("K136").Value = 4560
$ws.Range("M136").Value = -2010
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 82499.5
$ws.Range("H94").Value = 108624.55
$ws.Range("I94").Value = 145358
$ws.Range("J94").Value = 10668.667
$ws.Range("K94").Value = 145358
$ws.Range("L94").Value = 10668.667
$ws.Range("M94").Value = -144907
$ws.Range("N94").Value = -11570.667
$ws.Range("H105").Value = 2344.353
$ws.Range("I105").Value = 1803.4
$ws.Range("J105").Value = 3117.1428
$ws.Range("K105").Value = 1803.4
$ws.Range("L105").Value = 3117.1428
$ws.Range("M105").Value = -56.40000000000009
$ws.Range("N105").Value = -6611.1428
$ws.Range("H107").Value = 1520.3636
$ws.Range("I107").Value = 1433.1666
$ws.Range("J107").Value = 1625
$ws.Range("K107").Value = 1433.1666
$ws.Range("L107").Value = 1625
$ws.Range("M107").Value = 486.8334
$ws.Range("N107").Value = -5465
$ws.Range("H122").Value = 929.1429000000001
$ws.Range("I122").Value = 929.1429000000001
$ws.Range("K122").Value = 2787.4287
$ws.Range("M122").Value = -337.4287000000004
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 682.7143
$ws.Range("I12").Value = 658.6923
$ws.Range("J12").Value = 703.5333000000001
$ws.Range("K12").Value = 1976.0769
$ws.Range("L12").Value = 2110.5999
$ws.Range("M12").Value = -1803.0769
$ws.Range("N12").Value = -2456.5999
$ws.Range("H14").Value = 1041.6666
$ws.Range("I14").Value = 1041.6666
$ws.Range("K14").Value = 3124.9998
$ws.Range("M14").Value = -2951.9998
$ws.Range("H113").Value = 2109.25
$ws.Range("I113").Value = 1588.3334
$ws.Range("J113").Value = 2282.889
$ws.Range("K113").Value = 4765.0002
$ws.Range("L113").Value = 6848.667
$ws.Range("M113").Value = -2595.0002
$ws.Range("N113").Value = -11188.667
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6981.5
$ws.Range("I70").Value = 5983
$ws.Range("K70").Value = 5983
$ws.Range("M70").Value = -5713
$ws.Range("H73").Value = 6981.5
$ws.Range("I73").Value = 5983
$ws.Range("K73").Value = 5983
$ws.Range("M73").Value = -5047
$ws.Range("H94").Value = 39911.832
$ws.Range("J94").Value = 39911.832
$ws.Range("L94").Value = 39911.832
$ws.Range("N94").Value = -41263.832
$ws.Range("H113").Value = 4666.6665
$ws.Range("I113").Value = 1000
$ws.Range("K113").Value = 1000
$ws.Range("M113").Value = 1170
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 901
$ws.Range("I22").Value = 859
$ws.Range("K22").Value = 859
$ws.Range("M22").Value = -564
$ws.Range("H27").Value = 901
$ws.Range("I27").Value = 859
$ws.Range("K27").Value = 859
$ws.Range("M27").Value = -752
$ws.Range("H46").Value = 38480.93
$ws.Range("I46").Value = 85376.5
$ws.Range("K46").Value = 85376.5
$ws.Range("M46").Value = -85188.5
$ws.Range("H69").Value = 24999
$ws.Range("J69").Value = 24999
$ws.Range("L69").Value = 24999
$ws.Range("N69").Value = -26621
$ws.Range("H72").Value = 24999
$ws.Range("J72").Value = 24999
$ws.Range("L72").Value = 74997
$ws.Range("N72").Value = -83109
$ws.Range("H130").Value = 20000
$ws.Range("J130").Value = 20000
$ws.Range("L130").Value = 20000
$ws.Range("N130").Value = -30040
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 6500
$ws.Range("I52").Value = 3500
$ws.Range("K52").Value = 3500
$ws.Range("M52").Value = -3274
$ws.Range("H75").Value = 84624
$ws.Range("I75").Value = 84118
$ws.Range("K75").Value = 84118
$ws.Range("M75").Value = -83182
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = $null
$ws.Range("N76").Value = 0
$ws.Range("H78").Value = 84624
$ws.Range("I78").Value = 84118
$ws.Range("K78").Value = 252354
$ws.Range("M78").Value = -247674
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = $null
$ws.Range("N79").Value = 0
$ws.Range("H132").Value = 10058.857
$ws.Range("I132").Value = 8254.947
$ws.Range("K132").Value = 24764.841
$ws.Range("M132").Value = -22234.841
$ws.Range("H136").Value = 8916.546
$ws.Range("I136").Value = 8916.546
$ws.Range("K136").Value = 26749.638
$ws.Range("M136").Value = -24199.638
